$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the value of A17 to the new string
$ws.Range("A17").Value = "测试对比"

# Update the active selection to A17 (matches the recorded cursor move in the diff)
$ws.Range("A17").Select()
